# Update cryptos list values (price + 1h volume change) per upstream diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.868.69'
$ws.Range('E2').Value = '  +1.10%  '
$ws.Range('D3').Value = '1.841.80'
$ws.Range('E3').Value = '  +1.15%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.008'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '309.28'
$ws.Range('E5').Value = '  +1.35%  '
$ws.Range('E6').Value = '  -0.22%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4748'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3664'
$ws.Range('E8').Value = '  +2.06%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07192'
$ws.Range('E9').Value = '  +0.92%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9240'
$ws.Range('E10').Value = '  +3.10%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '19.63'
$ws.Range('E11').Value = '  +1.76%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07688'
$ws.Range('E12').Value = '  -1.00%  '
$ws.Range('D13').Value = '1.883.82'
$ws.Range('E13').Value = '  +3.29%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.303'
$ws.Range('E14').Value = '  +1.21%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.387'
$ws.Range('E15').Value = '  +1.16%  '
$ws.Range('E16').Value = '  +1.64%  '
$ws.Range('E17').Value = '  -0.15%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008617'
$ws.Range('E18').Value = '  +1.08%  '
$ws.Range('E19').Value = '  -0.15%  '
$ws.Range('D20').Value = '26.905.70'
$ws.Range('E20').Value = '  +1.06%  '
$ws.Range('E21').Value = '  +3.09%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.045'
$ws.Range('E22').Value = '  +0.74%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.62'
$ws.Range('E23').Value = '  +1.04%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.917'
$ws.Range('E24').Value = '  +0.23%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '152.24'
$ws.Range('E25').Value = '  +0.16%  '
$ws.Range('E26').Value = '  +1.55%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.995'
$ws.Range('E27').Value = '  +1.76%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '114.06'
$ws.Range('E28').Value = '  +0.53%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.918'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.08856'
$ws.Range('E30').Value = '  +0.73%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.304'
$ws.Range('E31').Value = '  +5.35%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.7461'
$ws.Range('E32').Value = '  +2.57%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.167'
$ws.Range('E33').Value = '  +3.82%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.472'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.717'
$ws.Range('E35').Value = '  -0.35%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.090'
$ws.Range('E36').Value = '  +1.74%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01948'
$ws.Range('E37').Value = '  +1.52%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.05254'
$ws.Range('E38').Value = '  +3.17%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.965'
$ws.Range('E39').Value = '  +1.72%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.5185'
$ws.Range('E40').Value = '  +3.16%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.934'
$ws.Range('E41').Value = '  +1.44%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1507'
$ws.Range('E42').Value = '  +1.14%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.166'
$ws.Range('E43').Value = '  +2.89%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '10.54'
$ws.Range('E44').Value = '  +5.63%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.4714'
$ws.Range('E45').Value = '  +1.78%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.006'
$ws.Range('E46').Value = '  -0.29%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '101.24'
$ws.Range('E47').Value = '  +3.20%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.595'
$ws.Range('E48').Value = '  +2.82%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '65.58'
$ws.Range('E49').Value = '  +3.29%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06022'
$ws.Range('E50').Value = '  +0.79%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.8851'
$ws.Range('E51').Value = '  +4.37%  '
